# "Test data for Greece Market"
#
# Before: Germany, Belgium, Czech, Swiss, Portugal, Slovakia, Italy, Spain, Croatia
# After : ... Croatia, Greece   (Greece added as a copy of Croatia, filled in)
#
# Along the way:
#   - Slovakia's lingering "select-all" window state is replaced with a
#     normal single-cell selection (F13).
#   - Croatia (which still held placeholder data copied from Slovakia) gets
#     its real market data filled in, and its window selection reverts to
#     the sheet's untouched/default state since focus moves to the new copy.
#   - The new sheet is renamed "Greece" and filled in with its own data,
#     ending up as the active/selected tab.

$wb = $excel.ActiveWorkbook

# Slovakia: give it a real selection instead of the "A1:XFD1048576" select-all
# state, without disturbing which sheet is actually active.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Range("F13").Select()

# Duplicate Croatia (last tab) to create the new country sheet right after it.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Cells.Select()
$croatia.Copy($null, $croatia)

$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Croatia previously carried over Slovakia's placeholder values - fill in
# its own data now.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T"

# Fill in the new Greece sheet's data (B4 entered before B2, matching the
# shared-string insertion order of the original edit).
$greece.Range("B4").Value = "NGC-4119/T3169"
$greece.Range("B2").Value = "Greece Market"

# Leave Greece as the active sheet/tab with B4 selected.
$greece.Range("B4").Select()
